# Regenerate the lattice-multiplication problems: each of the 15 table
# cells holds a 5-line block (problem, top factor digits, separator,
# two partial-product rows) joined by manual line breaks (Chr(11)).
# Replacing Cell.Range.Text rewrites the cell's run(s) in place while
# keeping the existing run formatting (sz=32) and the w:br line-break
# structure intact.
$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "82 x 93" + [char]11 + "  9    3" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "2|    |"
$tbl.Cell(1, 2).Range.Text = "90 x 53" + [char]11 + "  5    3" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "0|    |"
$tbl.Cell(1, 3).Range.Text = "13 x 22" + [char]11 + "  2    2" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "3|    |"
$tbl.Cell(2, 1).Range.Text = "69 x 43" + [char]11 + "  4    3" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "9|    |"
$tbl.Cell(2, 2).Range.Text = "90 x 49" + [char]11 + "  4    9" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "0|    |"
$tbl.Cell(2, 3).Range.Text = "66 x 47" + [char]11 + "  4    7" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "6|    |"
$tbl.Cell(3, 1).Range.Text = "78 x 22" + [char]11 + "  2    2" + [char]11 + "  ----" + [char]11 + "7|    |" + [char]11 + "8|    |"
$tbl.Cell(3, 2).Range.Text = "16 x 83" + [char]11 + "  8    3" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "6|    |"
$tbl.Cell(3, 3).Range.Text = "84 x 15" + [char]11 + "  1    5" + [char]11 + "  ----" + [char]11 + "8|    |" + [char]11 + "4|    |"
$tbl.Cell(4, 1).Range.Text = "96 x 20" + [char]11 + "  2    0" + [char]11 + "  ----" + [char]11 + "9|    |" + [char]11 + "6|    |"
$tbl.Cell(4, 2).Range.Text = "53 x 97" + [char]11 + "  9    7" + [char]11 + "  ----" + [char]11 + "5|    |" + [char]11 + "3|    |"
$tbl.Cell(4, 3).Range.Text = "33 x 94" + [char]11 + "  9    4" + [char]11 + "  ----" + [char]11 + "3|    |" + [char]11 + "3|    |"
$tbl.Cell(5, 1).Range.Text = "40 x 99" + [char]11 + "  9    9" + [char]11 + "  ----" + [char]11 + "4|    |" + [char]11 + "0|    |"
$tbl.Cell(5, 2).Range.Text = "62 x 72" + [char]11 + "  7    2" + [char]11 + "  ----" + [char]11 + "6|    |" + [char]11 + "2|    |"
$tbl.Cell(5, 3).Range.Text = "17 x 44" + [char]11 + "  4    4" + [char]11 + "  ----" + [char]11 + "1|    |" + [char]11 + "7|    |"
